$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new rows (22-30) following the existing pattern
$startRow = 22
for ($i = 0; $i -lt 9; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 10002 + $i
    $ws.Cells.Item($row, 2).Value = 110021 + $i
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Select the row after the data, matching the post-edit selection
$ws.Range("A31:XFD1048576").Select()

# Set up page setup similar to target (portrait orientation, 300 dpi print settings)
$ws.PageSetup.Orientation = 1
